$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / header row (values are unchanged text, but included for completeness)
$ws.Range("A1").Value = "SQA Signoff for conflict-times test files"

$ws.Range("A2").Value = "File Name"
$ws.Range("B2").Value = "SQA Name"
$ws.Range("C2").Value = "Date of Signoff"
$ws.Range("D2").Value = "Results"
$ws.Range("E2").Value = "Remarks"

# Rows 3-32: swap the class-times sign-off data for the conflict-times sign-off data.
# A column: " classt0XX" -> " conflictt0XX"
# B column: "Alla Salah" -> "Jared Cox"
# C column: 41385 (2013-04-21) -> 41386 (2013-04-22)
# D column ("pass"): kept for rows 3-24, removed (cleared) for rows 25-32
for ($i = 1; $i -le 30; $i++) {
    $r = $i + 2
    $fileNum = "{0:D3}" -f $i
    $ws.Cells.Item($r, 1).Value = " conflictt$fileNum"
    $ws.Cells.Item($r, 2).Value = "Jared Cox"
    $ws.Cells.Item($r, 3).Value = 41386

    if ($r -le 24) {
        $ws.Cells.Item($r, 4).Value = "pass"
    } else {
        $ws.Cells.Item($r, 4).Clear()
    }
}

# Update the sheet view: drop the stale scroll position and move the selection.
$ws.Range("D24").Select()
